$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '302.14'
Set-TextValue "E2" '-5.68%'

Set-TextValue "D3" '35.09'
Set-TextValue "E3" '-2.72%'

Set-TextValue "E4" '-1.63%'

Set-TextValue "D5" '0.07908'
Set-TextValue "E5" '-3.12%'

Set-TextValue "D6" '1.919'
Set-TextValue "E6" '-10.55%'

Set-TextValue "D7" '4.029'
Set-TextValue "E7" '-2.62%'

Set-TextValue "D8" '7.720'
Set-TextValue "E8" '-3.99%'

Set-TextValue "B9" 'MXToken'
Set-TextValue "C9" 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue "D9" '0.9236'
Set-TextValue "E9" '-0.27%'

Set-TextValue "B10" 'LiechtensteinCryptoassetsExchange'
Set-TextValue "C10" 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue "D10" '0.1205'
Set-TextValue "E10" '19.62%'

Set-TextValue "B11" 'WazirX'
Set-TextValue "C11" 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue "D11" '0.1846'
Set-TextValue "E11" '-2.28%'

Set-TextValue "B12" 'MandalaExchangeToken'
Set-TextValue "C12" 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue "D12" '0.09510'
Set-TextValue "E12" '2.67%'

Set-TextValue "B13" 'BitrueCoin'
Set-TextValue "C13" 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue "D13" '0.03538'
Set-TextValue "E13" '-1.47%'

Set-TextValue "B14" 'BitMartToken'
Set-TextValue "C14" 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue "D14" '0.09881'
Set-TextValue "E14" '-0.46%'

Set-TextValue "B15" 'BitForexToken'
Set-TextValue "C15" 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue "D15" '0.001403'
Set-TextValue "E15" '-2.23%'

Set-TextValue "B16" 'TigerCash'
Set-TextValue "C16" 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue "D16" '0.005833'
Set-TextValue "E16" '2.70%'

Set-TextValue "B17" 'LEO'
Set-TextValue "C17" 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue "D17" '3.491'
Set-TextValue "E17" '1.23%'

Set-TextValue "B18" 'BTSEToken'
Set-TextValue "C18" 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue "D18" '2.903'
Set-TextValue "E18" '3.67%'

Set-TextValue "E19" '2.11%'

Set-TextValue "D20" '0.1291'
Set-TextValue "E20" '-2.89%'

Set-TextValue "D21" '5.044'
Set-TextValue "E21" '-0.06%'

Set-TextValue "D22" '0.2399'
Set-TextValue "E22" '9.66%'

Set-TextValue "D23" '0.04496'
Set-TextValue "E23" '-1.92%'

Set-TextValue "E24" '-2.10%'

Set-TextValue "D25" '0.004570'
Set-TextValue "E25" '-3.38%'

Set-TextValue "E26" '-3.75%'

Set-TextValue "E27" '-6.79%'

Set-TextValue "D39" '0.01898'
Set-TextValue "E39" '-6.30%'

Set-TextValue "D40" '0.04716'
Set-TextValue "E40" '-5.65%'

Set-TextValue "D41" '0.007549'
Set-TextValue "E41" '-3.40%'

Set-TextValue "D42" '0.009554'
Set-TextValue "E42" '22.46%'

Set-TextValue "E43" '-5.50%'

Set-TextValue "E44" '0.78%'

Set-TextValue "D45" '0.01119'
Set-TextValue "E45" '-8.24%'

Set-TextValue "D46" '0.00006032'
Set-TextValue "E46" '-6.89%'

Set-TextValue "E47" '0.03%'

Set-TextValue "E49" '-31.33%'

Set-TextValue "E50" '0.03%'

Set-TextValue "E51" '0.03%'
